$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "MultColumnUnique"

# Headers - order matters for shared string table indices
$ws.Range("E4").Value = "Value"
$ws.Range("D4").Value = "Run"
$ws.Range("C4").Value = "Unit"
$ws.Range("B4").Value = "Current"

# Data rows
$data = @(
  @(15,1,1,1),
  @(15,2,1,2),
  @(15,3,1,3),
  @(15,4,1,4),
  @(20,1,2,3),
  @(20,2,2,3),
  @(20,3,2,4),
  @(20,4,2,4),
  @(26,1,3,4),
  @(26,2,3,4),
  @(26,3,3,4),
  @(26,4,3,4)
)
$r = 5
foreach ($row in $data) {
  $ws.Cells.Item($r, 2).Value = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
  $r++
}

$lo = $ws.ListObjects.Add(1, $ws.Range("B4:E16"), [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table5"

$ws.Range("H5").Formula2 = "=UNIQUE(Table5[[#All],[Current]:[Unit]])"
$ws.Range("J5").Value = "Mean"
$ws.Range("J6").Formula = "=AVERAGEIF(Table5[Current],H6,Table5[Value])"
$ws.Range("J7").Formula = "=AVERAGEIF(Table5[Current],H7,Table5[Value])"
$ws.Range("J8").Formula = "=AVERAGEIF(Table5[Current],H8,Table5[Value])"
$ws.Range("J9").Formula = "=AVERAGEIF(Table5[Current],H9,Table5[Value])"
$ws.Range("J10").Formula = "=AVERAGEIF(Table5[Current],H10,Table5[Value])"
$ws.Range("J11").Formula = "=AVERAGEIF(Table5[Current],H11,Table5[Value])"
$ws.Range("J12").Formula = "=AVERAGEIF(Table5[Current],H12,Table5[Value])"
$ws.Range("J13").Formula = "=AVERAGEIF(Table5[Current],H13,Table5[Value])"
$ws.Range("J14").Formula = "=AVERAGEIF(Table5[Current],H14,Table5[Value])"
$ws.Range("J15").Formula = "=AVERAGEIF(Table5[Current],H15,Table5[Value])"
$ws.Range("J16").Formula = "=AVERAGEIF(Table5[Current],H16,Table5[Value])"
$ws.Range("J17").Formula = "=AVERAGEIF(Table5[Current],H17,Table5[Value])"

$ws.Range("L4").Value = "The key is averageifs."
$ws.Range("L6").Formula2 = "=UNIQUE(CHOOSECOLS(Table5[],1,3))"
$ws.Range("N6").Formula = "=AVERAGEIFS(Table5[Value],Table5[Current],L6,Table5[Run],M6)"
$ws.Range("N7").Formula = "=AVERAGEIFS(Table5[Value],Table5[Current],L7,Table5[Run],M7)"
$ws.Range("N8").Formula = "=AVERAGEIFS(Table5[Value],Table5[Current],L8,Table5[Run],M8)"

$ws.Range("L5").Select()
